$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CTLT")

# Row 6: Change in inventories
$ws.Range("B6").Value = -205100000.0
$ws.Range("C6").Value = -142000000.0
$ws.Range("D6").Value = -76400000.0
$ws.Range("E6").Value = -45600000.0
$ws.Range("F6").Value = -19600000.0

# Row 8: Change in payables and accrued liability
$ws.Range("B8").Value = 1283000000.0
$ws.Range("C8").Value = 1169000000.0
$ws.Range("D8").Value = 787300000.0
$ws.Range("E8").Value = 481400000.0
$ws.Range("F8").Value = 204700000.0
